# Common: Added liquid volume
# Inserts a new "volume" column (C) on the "liquids" sheet, shifting the
# existing "pg" / "vg" / "description" columns one slot to the right
# (C->D, D->E, E->F), and fills in the new volume values per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("liquids")

# New volume values, keyed by row number.
$volumes = @{
    2  = 20
    3  = 20
    4  = 15
    5  = 10
    6  = 10
    7  = 10
    8  = 10
    9  = 20
    10 = 20
    11 = 20
    12 = 20
}

# Shift existing data (pg: C->D, vg: D->E, description: E->F) for each row,
# starting from the right-most column so we don't clobber data before it
# has been copied.
for ($r = 2; $r -le 12; $r++) {
    $pg = $ws.Cells.Item($r, 3).Value2
    $vg = $ws.Cells.Item($r, 4).Value2
    $desc = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 6).Value = $desc
    $ws.Cells.Item($r, 5).Value = $vg
    $ws.Cells.Item($r, 4).Value = $pg
    $ws.Cells.Item($r, 3).Value = $volumes[$r]
}

# Shift header row labels the same way, then set the new "volume" header.
$descHeader = $ws.Cells.Item(1, 5).Value2
$vgHeader = $ws.Cells.Item(1, 4).Value2
$pgHeader = $ws.Cells.Item(1, 3).Value2

$ws.Cells.Item(1, 6).Value = $descHeader
$ws.Cells.Item(1, 5).Value = $vgHeader
$ws.Cells.Item(1, 4).Value = $pgHeader
$ws.Cells.Item(1, 3).Value = "volume"

# Column widths: C/D keep their original widths (pg/vg data just shifted
# underneath them); E shrinks to fit "vg", and the new F column (the
# description) takes on the old E width-ish sizing.
$ws.Columns.Item(5).ColumnWidth = 8.17
$ws.Columns.Item(6).ColumnWidth = 25.5

# Move the selection to match the author's final cursor position.
[void]$ws.Range("B7").Select()
